# Update the division-practice answer table: each cell in the single
# 20-row/5-column table gets a new "a÷b=c, d" answer string. Addressing
# cells by (row, column) avoids any ambiguity from find/replace text
# collisions (several old/new answer strings repeat across the table).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "24÷3=8, 0"
$t.Cell(1,2).Range.Text = "55÷6=9, 1"
$t.Cell(1,3).Range.Text = "38÷4=9, 2"
$t.Cell(1,4).Range.Text = "88÷6=14, 4"
$t.Cell(1,5).Range.Text = "14÷5=2, 4"
$t.Cell(5,1).Range.Text = "75÷5=15, 0"
$t.Cell(5,2).Range.Text = "77÷7=11, 0"
$t.Cell(5,3).Range.Text = "70÷3=23, 1"
$t.Cell(5,4).Range.Text = "24÷9=2, 6"
$t.Cell(5,5).Range.Text = "52÷9=5, 7"
$t.Cell(9,1).Range.Text = "67÷3=22, 1"
$t.Cell(9,2).Range.Text = "73÷2=36, 1"
$t.Cell(9,3).Range.Text = "53÷4=13, 1"
$t.Cell(9,4).Range.Text = "48÷4=12, 0"
$t.Cell(9,5).Range.Text = "38÷9=4, 2"
$t.Cell(13,1).Range.Text = "19÷5=3, 4"
$t.Cell(13,2).Range.Text = "56÷4=14, 0"
$t.Cell(13,3).Range.Text = "11÷8=1, 3"
$t.Cell(13,4).Range.Text = "91÷2=45, 1"
$t.Cell(13,5).Range.Text = "67÷3=22, 1"
$t.Cell(17,1).Range.Text = "48÷3=16, 0"
$t.Cell(17,2).Range.Text = "99÷8=12, 3"
$t.Cell(17,3).Range.Text = "85÷9=9, 4"
$t.Cell(17,4).Range.Text = "26÷5=5, 1"
$t.Cell(17,5).Range.Text = "42÷3=14, 0"
